$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Occurrence 1 (first paragraph, italic run):
#   ", em favor de BMP MONEY PLUS SOCIEDADE DE CRÉDITO DIRETO S.A., ..."
#   -> ", em favor de BMP SOCIEDADE DE CRÉDITO DIRETO S.A., ..."
# The "_GoBack" bookmark ends up wrapping "BMP" after the edit, so we
# reposition it there BEFORE deleting " MONEY PLUS" -- this forces Word
# to split the run at that boundary instead of merging the whole
# sentence into one run.
# ---------------------------------------------------------------------
$full = $d.Content.Text
$idx = $full.IndexOf("BMP MONEY PLUS SOCIEDADE")

$rngBMP = $d.Range($idx, $idx + 3)
$d.Bookmarks.Add("_GoBack", $rngBMP)

$delStart = $idx + 3
$delEnd = $delStart + 11
$d.Range($delStart, $delEnd).Delete()

# ---------------------------------------------------------------------
# Occurrence 2 (numbered list paragraph):
#   "A BMP MONEY PLUS SOCIEDADE DE CRÉDITO DIRETO S.A., inscrita ..."
#   -> "A BMP SOCIEDADE DE CRÉDITO DIRETO S.A., inscrita ..."
# Use temporary bookmarks to pin the run-split boundaries around "BMP "
# and after the CNPJ number, matching how the run happens to be broken
# up in the edited document, then remove the scaffolding bookmarks.
# ---------------------------------------------------------------------
$full = $d.Content.Text
$idx2 = $full.IndexOf("BMP MONEY PLUS SOCIEDADE")

$rngA = $d.Range($idx2 - 2, $idx2)
$d.Bookmarks.Add("zzzTempA", $rngA)

$rngBMP2 = $d.Range($idx2, $idx2 + 4)
$d.Bookmarks.Add("zzzTempBMP", $rngBMP2)

$endMarker = "34.337.707/0001-00"
$mIdx = $full.IndexOf($endMarker, $idx2)
$endIdx = $mIdx + $endMarker.Length
$rngEnd = $d.Range($endIdx, $endIdx)
$d.Bookmarks.Add("zzzTempEnd", $rngEnd)

$delStart2 = $idx2 + 4
$delEnd2 = $delStart2 + 11
$d.Range($delStart2, $delEnd2).Delete()

$d.Bookmarks("zzzTempA").Delete()
$d.Bookmarks("zzzTempBMP").Delete()
$d.Bookmarks("zzzTempEnd").Delete()

# ---------------------------------------------------------------------
# Occurrence 3 (signature block, bold run, whole run is just the name):
#   "BMP MONEY PLUS SOCIEDADE DE CRÉDITO DIRETO S.A."
#   -> "BMP SOCIEDADE DE CRÉDITO DIRETO S.A."
# This run has no identically-formatted neighbour, so a plain Range.Text
# replace is safe here.
# ---------------------------------------------------------------------
$full = $d.Content.Text
$idx3 = $full.IndexOf("BMP MONEY PLUS SOCIEDADE DE CRÉDITO DIRETO S.A.")

$delStart3 = $idx3 + 3
$delEnd3 = $delStart3 + 11
$d.Range($delStart3, $delEnd3).Delete()
